$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C gets an explicit "applied" General number format (same visual format,
#     just marked as explicitly applied) and auto-fit width.
$ws.Range("C2:C11").NumberFormat = "General"

# --- Columns B & D get the "d/mm/yyyy;@" custom date format.
$ws.Range("B2:B11").NumberFormat = "d/mm/yyyy;@"
$ws.Range("D2:D11").NumberFormat = "d/mm/yyyy;@"

# --- Remove the erroring formula in D2, leave the cell blank (but keep its new format).
$ws.Range("D2").ClearContents()

# --- Fill in missing duration values in column C.
$ws.Range("C4").Value = 7
$ws.Range("C5").Value = 21
$ws.Range("C6").Value = 3
$ws.Range("C8").Value = 30

# --- Column C: best-fit (auto-fit) width, as seen after the user's edits.
$ws.Columns("C").AutoFit()

# --- Last user selection before save.
$ws.Range("B11").Select()
